$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a daily-frequency price table for "Femacal de La Calera -
# Frutilla". A new weekly observation is inserted as a new record right
# before the current row 298, pushing every row from 298 downward down by
# one (298->299, ..., 310->311).
$ws.Rows.Item(298).Insert()

# Fill in the newly-inserted row 298 with the new observation.
$ws.Cells.Item(298, 1).Value = 3
$ws.Cells.Item(298, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(298, 3).Value = "Coquimbo"
$ws.Cells.Item(298, 4).Value = 44706
$ws.Cells.Item(298, 5).Value = 5
$ws.Cells.Item(298, 6).Value = "Fruta"
$ws.Cells.Item(298, 7).Value = 100101
$ws.Cells.Item(298, 8).Value = "Berries"
$ws.Cells.Item(298, 9).Value = 100112025
$ws.Cells.Item(298, 10).Value = "Frutilla"
$ws.Cells.Item(298, 11).Value = "Sin especificar"
$ws.Cells.Item(298, 12).Value = "Primera"
$ws.Cells.Item(298, 13).Value = 56
$ws.Cells.Item(298, 14).Value = 7000
$ws.Cells.Item(298, 15).Value = 7000
$ws.Cells.Item(298, 16).Value = 7000
$ws.Cells.Item(298, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(298, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(298, 19).Value = 1000
$ws.Cells.Item(298, 20).Value = 7
